$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-12 from 45175 to 45183
$ws.Range("C2:C12").Value = 45183
